$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text (shared strings)
$ws.Range("B1").Value = "Creamy White"
$ws.Range("C1").Value = "Creamy Blue"

# Update data values (years + figures)
$ws.Range("A2").Value = 2016
$ws.Range("B2").Value = 30
$ws.Range("C2").Value = 47

$ws.Range("A3").Value = 2017
$ws.Range("B3").Value = 45
$ws.Range("C3").Value = 33

$ws.Range("A4").Value = 2018
$ws.Range("B4").Value = 64
$ws.Range("C4").Value = 66

$ws.Range("A5").Value = 2019
$ws.Range("B5").Value = 23
$ws.Range("C5").Value = 78

# Fix selection bug: select D1 instead of D2
$ws.Range("D1").Select()
